# "Name changed to SF for Singleton File"
# Adds a new "clear" label (red text) to A8, moves the active selection,
# and sets the sheet to portrait page orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labeled cell in row 8 (shares the style used in similar workbooks:
# a dedicated red font applied via a new cellXfs entry).
$ws.Range("A8").Value = "clear"
$ws.Range("A8").Font.Color = 255

# Move the active cell / selection as recorded by the workbook.
[void]$ws.Range("I15").Select()

# Page orientation -> portrait.
$ws.PageSetup.Orientation = 1
